# Ajustes para comparacion de grupo
# Update the "goals/points" column (B) for several teams in the
# "Participantes" sheet. Values are stored as text, so we prefix the
# new value with an apostrophe to force Excel to keep it as text
# instead of converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value  = "'6"
$ws.Cells.Item(2, 2).Value  = "'6"
$ws.Cells.Item(3, 2).Value  = "'3"
$ws.Cells.Item(4, 2).Value  = "'3"
$ws.Cells.Item(6, 2).Value  = "'3"
$ws.Cells.Item(7, 2).Value  = "'4"
$ws.Cells.Item(9, 2).Value  = "'6"
$ws.Cells.Item(10, 2).Value = "'0"
$ws.Cells.Item(11, 2).Value = "'9"
$ws.Cells.Item(12, 2).Value = "'3"
$ws.Cells.Item(13, 2).Value = "'4"
$ws.Cells.Item(15, 2).Value = "'5"
$ws.Cells.Item(16, 2).Value = "'4"
$ws.Cells.Item(17, 2).Value = "'1"
$ws.Cells.Item(18, 2).Value = "'1"
$ws.Cells.Item(19, 2).Value = "'6"
$ws.Cells.Item(20, 2).Value = "'9"
$ws.Cells.Item(23, 2).Value = "'6"
$ws.Cells.Item(24, 2).Value = "'0"
$ws.Cells.Item(25, 2).Value = "'6"
$ws.Cells.Item(26, 2).Value = "'6"
$ws.Cells.Item(27, 2).Value = "'3"
$ws.Cells.Item(28, 2).Value = "'3"
$ws.Cells.Item(29, 2).Value = "'1"
$ws.Cells.Item(30, 2).Value = "'6"
$ws.Cells.Item(32, 2).Value = "'6"
